$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap the "B quarter" and "C quarter" rows within each year block of 4
# rows (A,B,C,D quarters), e.g. row3<->row4, row7<->row8, ... Columns
# A,B,D,E are always swapped; column C is only touched when the two
# values actually differ, so cells that are blank on both sides are left
# completely untouched (avoids turning a "present but empty" cell into a
# fully-removed one where the data didn't actually change).
for ($base = 2; $base -le 81; $base += 4) {
    $rowB = $base + 1
    $rowC = $base + 2

    $rangeB1 = $ws.Range("A" + $rowB + ":B" + $rowB)
    $rangeC1 = $ws.Range("A" + $rowC + ":B" + $rowC)
    $valsB1 = $rangeB1.Value2
    $valsC1 = $rangeC1.Value2
    $rangeB1.Value2 = $valsC1
    $rangeC1.Value2 = $valsB1

    $rangeB2 = $ws.Range("D" + $rowB + ":E" + $rowB)
    $rangeC2 = $ws.Range("D" + $rowC + ":E" + $rowC)
    $valsB2 = $rangeB2.Value2
    $valsC2 = $rangeC2.Value2
    $rangeB2.Value2 = $valsC2
    $rangeC2.Value2 = $valsB2

    $cellB = $ws.Range("C" + $rowB)
    $cellC = $ws.Range("C" + $rowC)
    $valB = $cellB.Value2
    $valC = $cellC.Value2
    if ($valB -ne $valC) {
        $cellB.Value2 = $valC
        $cellC.Value2 = $valB
    }
}

# Drop the now-unwanted "产销率" / "销售量" (non-cumulative) columns F and G
$ws.Range("F:G").Delete()
